# Courts.xlsx - "Add files via upload / 31JanChanges"
#
# - Renames three sheets to remove spaces from their names.
# - Inserts two brand-new sheets ("CourtOfficerNotes" and "CourtReports")
#   right after "GenerateDocument", each with the usual
#   TESTCASE / SCRIPT_ITERATION / POM_ITERATION lead-in columns plus a
#   batch of new feature-specific columns, and a single stub data row.
# - Re-selects a couple of the untouched sheets (whole-row selections).
# - Leaves "CourtReports" as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing sheets (strip spaces from the display names).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Generate Document").Name   = "GenerateDocument"
$wb.Worksheets.Item("Generate Petition").Name   = "GeneratePetition"
$wb.Worksheets.Item("Petition Allegation").Name = "PetitionAllegation"

# ---------------------------------------------------------------------
# 2. Insert the two new worksheets right after "GenerateDocument".
#    Copy "CourtsApprovalAndAuditHistory" as a template sheet: it already
#    carries the shared TESTCASE/SCRIPT_ITERATION/POM_ITERATION header
#    styling (and no print/pageSetup baggage), then the template columns
#    get wiped and replaced with the new headers below.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("CourtsApprovalAndAuditHistory")
$anchor   = $wb.Worksheets.Item("GenerateDocument")

$template.Copy($null, $anchor)
$courtOfficerNotes = $wb.Worksheets.Item($anchor.Index + 1)
$courtOfficerNotes.Name = "CourtOfficerNotes"

$template.Copy($null, $courtOfficerNotes)
$courtReports = $wb.Worksheets.Item($courtOfficerNotes.Index + 1)
$courtReports.Name = "CourtReports"

# ---------------------------------------------------------------------
# 3. CourtOfficerNotes (A1:W2)
# ---------------------------------------------------------------------
$courtOfficerNotes.Range("D1:P2").Clear()

$courtOfficerNotes.Cells.Item(1,1).Value  = "TESTCASE"
$courtOfficerNotes.Cells.Item(1,2).Value  = "SCRIPT_ITERATION"
$courtOfficerNotes.Cells.Item(1,3).Value  = "POM_ITERATION"
$courtOfficerNotes.Cells.Item(1,4).Value  = "PARTICIPANTS_PRESENT_IN_COURT"
$courtOfficerNotes.Cells.Item(1,5).Value  = "PARTICIPANTS_PRESENT_REMOTELY"
$courtOfficerNotes.Cells.Item(1,6).Value  = "TAFF_WORKER_PRESENT_IN_COURT"
$courtOfficerNotes.Cells.Item(1,7).Value  = "STAFF_WORKER_PRESENT_REMOTELY"
$courtOfficerNotes.Cells.Item(1,8).Value  = "OTHER_PARTICIPANTS_ATTENDANCE"
$courtOfficerNotes.Cells.Item(1,9).Value  = "IMMEDIATE_ACTION_FORTH_WITH_ORDER_DUE_DATE"
$courtOfficerNotes.Cells.Item(1,10).Value = "IMMEDIATE_ACTION_FORTHWITH_ORDERS_ORDERED_BY_COURT"
$courtOfficerNotes.Cells.Item(1,11).Value = "DISCUSSION"
$courtOfficerNotes.Cells.Item(1,12).Value = "PLACEMENT_DISCUSSION"
$courtOfficerNotes.Cells.Item(1,13).Value = "ORDERS"
$courtOfficerNotes.Cells.Item(1,14).Value = "VISITATION"
$courtOfficerNotes.Cells.Item(1,15).Value = "NOTES"
$courtOfficerNotes.Cells.Item(1,16).Value = "IMMEDIATE_ACTION_FORTHWITH_ORDERS_ORDERED_BY_COURT"
$courtOfficerNotes.Cells.Item(1,17).Value = "CASA_APPOINTED_DETAILS"
$courtOfficerNotes.Cells.Item(1,18).Value = "ORDER_TO_SHOW_CAUSE"
$courtOfficerNotes.Cells.Item(1,19).Value = "DISCOVERY_DUE_DATE"
$courtOfficerNotes.Cells.Item(1,20).Value = "NEXT_HEARING_TYPE"
$courtOfficerNotes.Cells.Item(1,21).Value = "DEPARTMENT_FOR_NEXT_HEARING"
$courtOfficerNotes.Cells.Item(1,22).Value = "ROOM_FOR_NEXT_HEARING"
$courtOfficerNotes.Cells.Item(1,23).Value = "SAVE_BTN"

$courtOfficerNotes.Columns.Item(1).AutoFit()
$courtOfficerNotes.Columns.Item(2).AutoFit()
$courtOfficerNotes.Columns.Item(3).AutoFit()
for ($c = 4; $c -le 22; $c++) {
    $courtOfficerNotes.Columns.Item($c).AutoFit()
}

$courtOfficerNotes.Range("W1").Select()

# ---------------------------------------------------------------------
# 4. CourtReports (A1:T2)
# ---------------------------------------------------------------------
$courtReports.Range("D1:P2").Clear()

$courtReports.Cells.Item(1,1).Value  = "TESTCASE"
$courtReports.Cells.Item(1,2).Value  = "SCRIPT_ITERATION"
$courtReports.Cells.Item(1,3).Value  = "POM_ITERATION"
$courtReports.Cells.Item(1,4).Value  = "REPORT_TYPE"
$courtReports.Cells.Item(1,5).Value  = "STATUS"
$courtReports.Cells.Item(1,6).Value  = "DUE_DATE_TO_SUPERVISOR"
$courtReports.Cells.Item(1,7).Value  = "DUE_DATE_TO_CLERICAL"
$courtReports.Cells.Item(1,8).Value  = "DESCRIPTION"
$courtReports.Cells.Item(1,9).Value  = "TRANSLATION_DOCUMENTS_SENT_TO"
$courtReports.Cells.Item(1,10).Value = "METHOD_OF_TRANSLATION_DOCUMENTS_SENT"
$courtReports.Cells.Item(1,11).Value = "TRANSLATION_DOCUMENTS_SENT_DATE"
$courtReports.Cells.Item(1,12).Value = "TRANSLATION_DOCUMENTS_RECEIVED_DATE"
$courtReports.Cells.Item(1,13).Value = "NAME_OF_COURT"
$courtReports.Cells.Item(1,14).Value = "METHOD_OF_FILING"
$courtReports.Cells.Item(1,15).Value = "DATE_REPORT_FILED"
$courtReports.Cells.Item(1,16).Value = "FILING_DUE_DATE"
$courtReports.Cells.Item(1,17).Value = "APPROVAL_SUPERVISOR"
$courtReports.Cells.Item(1,18).Value = "SUBMITTED_FOR_APPROVAL_DATE"
$courtReports.Cells.Item(1,19).Value = "SUPERVISOR_APPROVAL_DATE"
$courtReports.Cells.Item(1,20).Value = "SAVE_BTN"

$courtReports.Columns.Item(1).AutoFit()
$courtReports.Columns.Item(2).AutoFit()
$courtReports.Columns.Item(3).AutoFit()
for ($c = 4; $c -le 19; $c++) {
    $courtReports.Columns.Item($c).AutoFit()
}

$courtReports.Range("F13").Select()

# ---------------------------------------------------------------------
# 5. Tidy up selections on a couple of the untouched sheets.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("GeneratePetition").Rows.Item(1).Select()
$wb.Worksheets.Item("CourtsApprovalAndAuditHistory").Range("A1:XFD2").Select()

# ---------------------------------------------------------------------
# 6. Leave "CourtReports" as the active tab.
# ---------------------------------------------------------------------
$courtReports.Activate()
$courtReports.Range("F13").Select()
